# 7.10 Fixed Some Bugs
# - Wrap the three "wrong choice" hint strings in green color-tag markup.
# - Grow rows 2 & 3 to fit the now two-line hint text.
# - Move the active selection to B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = " <color=#00CC00>(Don’t fall into habitual thinking.)</color>"
$ws.Range("B3").Value2 = " <color=#00CC00>(Carefully recall the investigation process.)</color>"
$ws.Range("B4").Value2 = " <color=#00CC00>(Where was the scene that clearly showed signs of a struggle?)</color>"

$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34

$ws.Range("B10").Select() | Out-Null
